$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.047.25"
$ws.Range("E2").Value = "  +3.89%  "
$ws.Range("D3").Value = "2.421.21"
$ws.Range("E3").Value = "  +3.41%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.360"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.07%  "
$ws.Range("D14").Value = "2.846.60"
$ws.Range("E14").Value = "  +3.22%  "
$ws.Range("D15").Value = "59.950.71"
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000139"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").Value = "2.415.55"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "331.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.171"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "0.0₃0783"
$ws.Range("E28").Value = "  +6.78%  "
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("B32").Value = "SuiNetwork"
$ws.Range("C32").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.25%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.56%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  +10.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "319.43"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0961"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.63%  "
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.29%  "
$ws.Range("B47").Value = "Polygon"
$ws.Range("C47").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.409"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.11%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.574"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0227"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("E51").Value = "  -0.22%  "
